$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "60.232.81"
$ws.Range("E2").Value = "  -6.09%  "
Set-TextCell $ws.Range("D3") "3.291.31"
$ws.Range("E3").Value = "  -5.23%  "
Set-TextCell $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextCell $ws.Range("D5") "558.89"
$ws.Range("E5").Value = "  -4.32%  "
Set-TextCell $ws.Range("D6") "127.61"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("E7").Value = "  -0.37%  "
Set-TextCell $ws.Range("D8") "3.289.88"
$ws.Range("E8").Value = "  -5.20%  "
Set-TextCell $ws.Range("D9") "0.473"
$ws.Range("E9").Value = "  -1.79%  "
Set-TextCell $ws.Range("D10") "7.32"
$ws.Range("E10").Value = "  -3.81%  "
$ws.Range("E11").Value = "  -4.92%  "
$ws.Range("E12").Value = "  -4.34%  "
Set-TextCell $ws.Range("D13") "3.847.96"
$ws.Range("E13").Value = "  -5.40%  "
Set-TextCell $ws.Range("D14") "0.119"
$ws.Range("E14").Value = "  -0.11%  "
Set-TextCell $ws.Range("D15") "3.284.72"
$ws.Range("E15").Value = "  -5.46%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell $ws.Range("D17") "60.399.26"
$ws.Range("E17").Value = "  -5.78%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell $ws.Range("D18") "24.11"
$ws.Range("E18").Value = "  -0.62%  "
Set-TextCell $ws.Range("D19") "5.63"
$ws.Range("E19").Value = "  -0.93%  "
Set-TextCell $ws.Range("D20") "13.28"
$ws.Range("E20").Value = "  -1.00%  "
Set-TextCell $ws.Range("D21") "9.05"
$ws.Range("E21").Value = "  -9.24%  "
Set-TextCell $ws.Range("D22") "352.00"
$ws.Range("E22").Value = "  -8.41%  "
Set-TextCell $ws.Range("D23") "0.551"
$ws.Range("E23").Value = "  -2.85%  "
Set-TextCell $ws.Range("D24") "1.00"
$ws.Range("E24").Value = "  +0.00%  "
Set-TextCell $ws.Range("D25") "3.417.07"
$ws.Range("E25").Value = "  -5.45%  "
Set-TextCell $ws.Range("D26") "69.33"
$ws.Range("E26").Value = "  -7.52%  "
Set-TextCell $ws.Range("D27") "0.0000108"
$ws.Range("E27").Value = "  -2.81%  "
Set-TextCell $ws.Range("D28") "0.995"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("E29").Value = "  +1.88%  "
$ws.Range("E30").Value = "  -0.47%  "
Set-TextCell $ws.Range("D31") "7.78"
$ws.Range("E31").Value = "  -1.82%  "
Set-TextCell $ws.Range("D32") "2.09"
$ws.Range("E32").Value = "  -5.98%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -1.93%  "
Set-TextCell $ws.Range("D35") "3.319.20"
$ws.Range("E35").Value = "  -5.21%  "
Set-TextCell $ws.Range("D36") "22.64"
$ws.Range("E36").Value = "  -0.97%  "
Set-TextCell $ws.Range("D37") "5.24"
$ws.Range("E37").Value = "  +1.11%  "
Set-TextCell $ws.Range("D38") "6.72"
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  -1.12%  "
Set-TextCell $ws.Range("D40") "158.04"
$ws.Range("E40").Value = "  -2.60%  "
Set-TextCell $ws.Range("D41") "0.0746"
$ws.Range("E41").Value = "  -3.97%  "
Set-TextCell $ws.Range("D42") "1.00"
$ws.Range("E42").Value = "  +0.00%  "
Set-TextCell $ws.Range("D43") "41.08"
$ws.Range("E43").Value = "  -0.59%  "
Set-TextCell $ws.Range("D44") "4.35"
$ws.Range("E44").Value = "  +1.53%  "
Set-TextCell $ws.Range("D45") "0.738"
$ws.Range("E45").Value = "  -7.32%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D47") "22.63"
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell $ws.Range("D48") "1.54"
$ws.Range("E48").Value = "  -4.88%  "
Set-TextCell $ws.Range("D49") "6.64"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E50").Value = "  -5.12%  "
Set-TextCell $ws.Range("D51") "21.19"
$ws.Range("E51").Value = "  +3.69%  "
